# Day 7 routing update:
# - Finish the "also accepts a regex" sentence (merge the two runs split by
#   the _GoBack bookmark) and extend the app.use discussion with the new
#   "did ... in previous versions, but now accepts regex." wording.
# - Insert a page break + "DAY 7 ROUTING" section heading + a new
#   "Routing and MVC architecture of REST API" line after it.
$d = $word.ActiveDocument

# Locate the target paragraph ("app.all also accepts a re...") via Find,
# then resolve it to the real enclosing Paragraph so we operate on the
# whole paragraph range (Find's own Range can end up scoped to just the
# matched text, not the full paragraph).
$searchRng = $d.Content
$found = $searchRng.Find.Execute("app.all also accepts a re", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find target paragraph (app.all also accepts a re...)" }

$paragraphs = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $candidate = $paragraphs.Item($i)
    if ($candidate.Range.Start -le $searchRng.Start -and $candidate.Range.End -ge $searchRng.End) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) { throw "Could not resolve enclosing paragraph for match" }

$targetRange = $paragraphs.Item($targetIndex).Range

# Replace the whole paragraph with the rewritten paragraph plus the four
# new paragraphs (page break, heading, blank spacer, and the new line) that
# follow it -- InsertXML on a full paragraph Range replaces its contents.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00175CC9" w:rsidRPr="00EB264A" w:rsidRDefault="00175CC9" w:rsidP="00175CC9"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:spacing w:before="240"/><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00175CC9"><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>app.all</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00175CC9"><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> also accepts a regex as its path parameter. </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>app.use</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> d</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>id</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> not accept a regex</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> in previous versions,</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>, but</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> now accepts regex</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:before="240"/><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:lastRenderedPageBreak/><w:t>DAY 7 ROUTING</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:before="240"/><w:jc w:val="both"/><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:before="240"/><w:jc w:val="center"/><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Routing and MVC architecture</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="00175CC9"><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> of REST API</w:t></w:r></w:p>'
$targetRange.InsertXML($xml)
